# repull data, push all data, mean calculation
# Update column F (dSF) values for the rows whose underlying data was
# repulled / recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = -1
    8  = -7
    9  = -2
    10 = 4
    11 = -7
    14 = -8
    15 = -4
    16 = -3
    17 = -5
    19 = -9
    20 = -7
    23 = -8
    24 = -3
    25 = -3
    26 = -4
    27 = -6
    31 = 0
    33 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
